$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Files": cells C3 and C5 each have their blank-line separators
# (encoded in the OOXML as runs of four "_x000d_" escapes, i.e. four CR
# characters, followed by a line break) extended by one extra CR, so the
# runs become five CRs instead of four.
#
# NOTE: this host's Range.Value getter normalises CR/CRLF to LF when
# *reading* a cell, so we can't reliably detect the existing CR runs by
# reading back the current value. Instead we rebuild the exact target
# strings from their known text segments, joining them with the new
# (five-CR + LF) separator directly.
# ---------------------------------------------------------------------
$wsFiles = $wb.Worksheets.Item("Files")

$cr = [string][char]13
$lf = [string][char]10
$sep = "$cr$cr$cr$cr$cr$lf"

$c3Segments = @(
  'Baseline characteristics of the 20 nursing homes (Table 1).',
  '',
  'N.B. Table 1 Questionnaire on Palliative care for Advanced Dementia (qPAD) median scores by group calculated using staff level data in Dataset 2.'
)
$wsFiles.Range("C3").Value = ($c3Segments -join $sep)

$c5Segments = @(
  'Contains ',
  '',
  'baseline characteristics of the 131 residents (Table 1)',
  '',
  'End of Life Outcomes (family-rated and nurse-rated) total scores and weeks to completion (Tables 2, 3 and 4)',
  '',
  'N.B. there are more ‘time to completion’ data points than EOLD (End Of Life Dementia) scores since some EOLD responses were provided but unable to generate a valid total score.',
  '',
  'Case conference (whether or not resident received at least one case conference)',
  '',
  'Quality of Life in Late-stage Dementia (QUALID) data',
  '',
  'Symptoms and care during the last month of life (Table 5)'
)
$wsFiles.Range("C5").Value = ($c5Segments -join $sep)

# ---------------------------------------------------------------------
# Sheet "People": remove the contactPoint columns (J1 header / J3 value)
# that get moved over to the "Organisations" sheet below.
# ---------------------------------------------------------------------
$wsPeople = $wb.Worksheets.Item("People")
$wsPeople.Range("J1").ClearContents()
$wsPeople.Range("J3").ClearContents()

# ---------------------------------------------------------------------
# Sheet "Organisations": add the contactPoint columns - G1 header, plus
# the contact info for Tim Luckett's organisation on rows 3 and 16.
# ---------------------------------------------------------------------
$wsOrg = $wb.Worksheets.Item("Organisations")
$wsOrg.Range("G1").Value = "contactPoint>TYPE:ContactPoint>"
$wsOrg.Range("G3").Value = "ID: tim.luckett@uts.edu.au, contactType: customer service, email: tim.luckett@uts.edu.au, name: Contact Tim Luckett"
$wsOrg.Range("G16").Value = "ID: tim.luckett@uts.edu.au, contactType: customer service, email: tim.luckett@uts.edu.au, name: Contact Tim Luckett"
